$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Range("G2")

# Replace the cell text with the updated, numbered note text.
$cell.Value = "Task for Commit 3: 1. The code to setup and launch WebDriver (Firefox) should be improved to include Location   Firefox.exe    and a Firefox Profile must be created.`n2. The amount of time the driver should wait when searching for a GUI element should be specified`n3.The code to setup and launch Firefox should be created as a separate method.`n4. All script initialization parameters like Location of Firefox , Base URL, User Name , Password etc should be stored in separate file say Util.java.  This helps in easy maintenance of script"

# Re-apply rich-text (per-character) bold formatting matching the
# original "Task for Commit 3" label plus the new "N." lead-ins,
# while the rest of each sentence stays regular (non-bold).
$cell.Characters(1, 17).Font.Bold = $true
$cell.Characters(18, 2).Font.Bold = $false
$cell.Characters(20, 25).Font.Bold = $true
$cell.Characters(45, 122).Font.Bold = $false
$cell.Characters(167, 21).Font.Bold = $true
$cell.Characters(188, 77).Font.Bold = $false
$cell.Characters(265, 24).Font.Bold = $true
$cell.Characters(289, 55).Font.Bold = $false
$cell.Characters(344, 42).Font.Bold = $true
$cell.Characters(386, 150).Font.Bold = $false
